# "updated title txt on slide 5"
# Slide 5's title currently reads "Create the Gui" (as two runs: "Create the "
# and "Gui"). Re-split and retype it as "Create " + "the Menu" so the title
# now reads "Create the Menu".

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(5)
$title = $slide.Shapes.Item(1)
$tr = $title.TextFrame.TextRange

$firstRun = $tr.Characters(1, 7)
$firstRun.Text = "Create "

$secondRun = $tr.Characters(8, 8)
$secondRun.Text = "the Menu"
